$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$tf = $s.Shapes.Item(3).TextFrame
$tr = $tf.TextRange

# --- Paragraph: "In our experiment, we chose the pytorch library ... heart failure." ---
$oldPara1 = "In our experiment, we chose the pytorch library to build a neural network for this classification task. The motivation was to use an available deep learning library that provides various optmizers and loss functions and an easy way to tune hyperparamters so we can arrive at the most accurate network. In this experiment we found that a neural network with 2 layers performed best with our data set of heart data. we ran the experiment with 8 -14 neurons in the first layer and 2 neurons in the final layer with the leakyrelu activation to predict the chances of heart failure."
$newPara1 = "In our experiment, we chose the PyTorch library to build a neural network for this classification task. The motivation was to use an available deep learning library that provides various optimizers and loss functions and an easy way to tune hyperparameters so we can arrive at the most accurate network. In this experiment we found that a neural network with 2 layers performed best with our data set of heart data. we ran the experiment with 8 -14 neurons in the first layer and 2 neurons in the final layer with the LeakyReLU activation to predict the chances of heart failure."

$full = $tr.Text
$idx1 = $full.IndexOf($oldPara1)
$start1 = $idx1 + 1
$len1 = $oldPara1.Length
$range1 = $tr.Characters($start1, $len1)
$range1.Text = $newPara1

# --- Paragraph: "We were able to achieve a top accuracy of 86% ... RMSProp." ---
$oldPara2 = "We were able to achieve a top accuracy of 86% over the test data with 12 neurons in the first layer and with optimizer set as RMSProp."
$newPara2 = "We were able to achieve a top accuracy of 86% over the test data with 12 neurons in the first layer and with optimizer set as RMSProp."

$full2 = $tr.Text
$idx2 = $full2.IndexOf($oldPara2)
$start2 = $idx2 + 1
$len2 = $oldPara2.Length
$range2 = $tr.Characters($start2, $len2)
$range2.Text = $newPara2
